$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("diagnostics")

$ws.Range("H6").Value = "c633b3cf-d150-4007-88ac-d18907a85098"
$ws.Range("H7").Value = "4d025129-d1ba-4cbb-ba9c-7720d2bc129c"
$ws.Range("H8").Value = "39b8f6a9-937b-4166-bf55-d361e794afd7"
$ws.Range("H9").Value = "842f317d-f770-4665-9316-dc84e72a5bca"
$ws.Range("H10").Value = "2d26f2f3-b9d5-408e-8558-7443322ecb7c"
$ws.Range("H11").Value = "e9085274-1009-47b2-9031-6bbd6a7f7377"
$ws.Range("H12").Value = "464b6b4e-839a-4bf0-b5a8-a77288e85245"
$ws.Range("H13").Value = "2dadf332-c4d5-41fd-9192-68715d02610a"
$ws.Range("H14").Value = "26a9f2ac-65c2-4f39-af08-1b26f402284f"
$ws.Range("H15").Value = "6e5c8133-7ae8-41e5-9edd-18e2c744288a"
$ws.Range("H16").Value = "3c080f0c-5745-4d05-b01b-0b4c4cd5ab9a"
$ws.Range("H17").Value = "75f94655-9307-4ca3-aabb-8fce81454578"
$ws.Range("H19").Value = "40169ffd-6050-45e5-93a8-caed0c5a1ab7"
$ws.Range("H20").Value = "b286ce66-a092-4d45-a87b-fb1987d1b7fa"
$ws.Range("H21").Value = "6aaed81c-8cf6-428f-a675-1399b03b6d82"
$ws.Range("H22").Value = "ee155508-eadc-4039-a375-c17b671833d7"
$ws.Range("H23").Value = "a6c775a1-d7d1-4ad4-924b-ef5ee0d1ef77"
$ws.Range("H24").Value = "e659eb16-a838-440e-95f2-e824db80bef9"
$ws.Range("H25").Value = "f0bfb036-a45b-4573-bb2b-755504461ec2"
$ws.Range("H26").Value = "27e1948a-4d16-48a3-a5b6-feffe8301427"
$ws.Range("H27").Value = "17faa7ec-5e65-48e6-9a41-f0819e681784"
$ws.Range("H28").Value = "d8032654-96d3-4656-8472-8a07a099df9f"
$ws.Range("H29").Value = "f55d501a-e9dc-42ab-b77c-ed29b6c6c0ef"
$ws.Range("H30").Value = "6ee7c933-3b81-4dd7-88f7-33f3e837f737"
$ws.Range("H31").Value = "2d2eae99-8fb3-446d-a2f8-e733d637e0d8"
$ws.Range("H32").Value = "6ce1dd33-2777-40f5-949c-0043cd21f1dc"
$ws.Range("H33").Value = "b30078a9-76a6-493b-96b0-7de65272fdc2"
$ws.Range("H34").Value = "a51f96ac-4e2a-41b7-a3f5-2ac184a421d1"
$ws.Range("H35").Value = "bcb2899d-7b3b-4aed-8536-8e24aecadc6e"
$ws.Range("H36").Value = "61b9d4fe-cec1-4171-9545-3adad30c8c41"
$ws.Range("H37").Value = "c0dd763c-f4bb-4d5b-93e0-e720bb183896"
$ws.Range("H38").Value = "a97fb3f2-382f-4c67-ac17-92a414689288"
$ws.Range("H39").Value = "e7d8c6bd-f14a-4882-846c-c8ceae6204fe"
$ws.Range("H40").Value = "3a3dc980-8f95-4c78-ad5a-1d315d45e717"
$ws.Range("H41").Value = "7377073d-93f7-42ce-b529-985b970397ce"
$ws.Range("H42").Value = "a7f7818a-771c-4c1d-bc1a-512526b5a7ae"
$ws.Range("H43").Value = "fe240662-989e-4984-aee0-1b8b83305bcb"
$ws.Range("H44").Value = "fc009024-5c7e-4aa2-b583-8f18428a6be3"
$ws.Range("H45").Value = "eccb1c32-c22b-4a8b-b699-01e01402d31b"
$ws.Range("H46").Value = "7a6315d3-b727-4a2c-b6fc-9862024579b2"
$ws.Range("H47").Value = "24200d78-47c6-41a5-b535-b77dd403e1ab"
$ws.Range("H48").Value = "a77709f3-1888-430e-a108-6ab7320065a8"
$ws.Range("H49").Value = "dcd52a6b-0684-4534-a8da-e262ab510d61"
$ws.Range("H50").Value = "bbe68dc1-da98-4e00-844b-e6e851a800e2"
$ws.Range("H51").Value = "8cd46abd-96ed-41da-90b2-d296365f18a9"
$ws.Range("H52").Value = "75177149-3acf-44f9-bb6f-2804bdd02f63"
$ws.Range("H53").Value = "ce84a5dc-1458-474b-a883-43df2e90e05b"
$ws.Range("H54").Value = "50afc1e9-8994-4006-9fb2-60b4ff258514"
$ws.Range("H55").Value = "cafe744b-c5fb-4561-b52c-3ab6002a05f5"
$ws.Range("H56").Value = "6f95f12f-fd70-4d08-b5b7-befe7de06e0f"
$ws.Range("H57").Value = "87d37215-2d02-4992-8725-1e5037655f97"
$ws.Range("H58").Value = "82bdaedf-d4e4-4639-b3f1-095e65770c5c"
$ws.Range("H59").Value = "e296eaa3-4bb2-474d-b797-114f43e83b59"
$ws.Range("H60").Value = "258ef7f3-0676-4df9-8770-c8715c3b0da7"
$ws.Range("H61").Value = "f581fb7b-4c62-4492-aa20-ce0e634754d2"
$ws.Range("H62").Value = "b30db91d-38c5-403d-8eba-d48c4e9b621c"
$ws.Range("H63").Value = "0519b8bf-ab6c-402e-81da-a60c36868381"
$ws.Range("H64").Value = "d7895669-c69b-445d-a381-cbdf6e463222"
$ws.Range("H65").Value = "d9d36687-6ace-4251-9007-f3d6693c23d5"
$ws.Range("H66").Value = "9b5db90f-411d-46b1-bfca-90e8887c1c3c"
